# Append three new drink rows (коффе/coffee items) below the existing
# water/cola/chips/juice/sandwich list, matching the style already used
# by rows 2-5 (column A mirrors the "0.000" numeric style already (mis)used
# for the text column, column B uses the integer style).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = "капучино"
$ws.Range("B6").Value = 40
$ws.Range("A7").Value = "латте"
$ws.Range("B7").Value = 40
$ws.Range("A8").Value = "раф"
$ws.Range("B8").Value = 60

# Match the formatting already applied to rows 2-5 (A: "0.000", B: "0"),
# so no new style entries get created — the new rows reuse the existing
# cellXfs (s="2" for A, s="3" for B).
$ws.Range("A6:A8").NumberFormat = "0.000"
$ws.Range("B6:B8").NumberFormat = "0"
$ws.Range("A6:B8").RowHeight = 30

# Leave the selection on the last-entered cell, like the source file.
[void]$ws.Range("B8").Select()
